$d = $word.ActiveDocument

# --- Step 1: Replace the lone "2. " paragraph with the full "Erro ao imprimir o ticket"
#     section (heading + three body paragraphs + embedded code screenshot textbox) ---
$p = $d.Paragraphs.Item(15)
$r = $p.Range
$xml = '<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" mc:Ignorable="w14">
<w:body>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Body"/>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>2.  Erro ao imprimir o ticket</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Body"/>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>Segundo especifica</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>çõ</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">es do sistema ao ser impresso um ticket </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">é </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>necess</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>á</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">rio verificar se o valor inserido </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">é </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>o suficiente para a impress</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>ã</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>o do ticket e ap</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>ó</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>s a impress</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>ã</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">o </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">é </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>necess</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>á</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>rio realizar o debito do valor do ticket do dinheiro inserido pelo cliente.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Body"/>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>no c</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>ó</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">digo original </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">é </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>realizada a verifica</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>çã</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>o do saldo do cliente mas n</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>ã</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">o </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">é </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>realizado o debito ap</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>ó</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>s a impress</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>ã</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>o do ticket. Como podemos observar abaixo:</w:t>
      </w:r>
      <w:r>
        <mc:AlternateContent>
          <mc:Choice Requires="wps">
            <w:drawing>
              <wp:anchor distT="152400" distB="152400" distL="152400" distR="152400" simplePos="0" relativeHeight="251660288" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
                <wp:simplePos x="0" y="0"/>
                <wp:positionH relativeFrom="margin">
                  <wp:posOffset>-6349</wp:posOffset>
                </wp:positionH>
                <wp:positionV relativeFrom="line">
                  <wp:posOffset>340439</wp:posOffset>
                </wp:positionV>
                <wp:extent cx="5727699" cy="2514561"/>
                <wp:effectExtent l="0" t="0" r="0" b="0"/>
                <wp:wrapTopAndBottom distT="152400" distB="152400"/>
                <wp:docPr id="1073741826" name="officeArt object"/>
                <wp:cNvGraphicFramePr/>
                <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
                  <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
                    <wps:wsp>
                      <wps:cNvSpPr/>
                      <wps:spPr>
                        <a:xfrm>
                          <a:off x="0" y="0"/>
                          <a:ext cx="5727699" cy="2514561"/>
                        </a:xfrm>
                        <a:prstGeom prst="rect">
                          <a:avLst/>
                        </a:prstGeom>
                        <a:solidFill>
                          <a:srgbClr val="000000"/>
                        </a:solidFill>
                        <a:ln w="12700" cap="flat">
                          <a:noFill/>
                          <a:miter lim="400000"/>
                        </a:ln>
                        <a:effectLst/>
                      </wps:spPr>
                      <wps:txbx>
                        <w:txbxContent>
                          <w:p>
                            <w:pPr>
                              <w:pStyle w:val="Body"/>
                              <w:jc w:val="left"/>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                              </w:rPr>
                            </w:pPr>
                            <w:r>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                                <w:rtl w:val="0"/>
                                <w:lang w:val="en-US"/>
                              </w:rPr>
                              <w:t xml:space="preserve"> public String imprimir() throws SaldoInsuficienteException {</w:t>
                            </w:r>
                          </w:p>
                          <w:p>
                            <w:pPr>
                              <w:pStyle w:val="Body"/>
                              <w:jc w:val="left"/>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                              </w:rPr>
                            </w:pPr>
                            <w:r>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                                <w:rtl w:val="0"/>
                              </w:rPr>
                              <w:t xml:space="preserve">        if (saldo &lt; valor) {</w:t>
                            </w:r>
                          </w:p>
                          <w:p>
                            <w:pPr>
                              <w:pStyle w:val="Body"/>
                              <w:jc w:val="left"/>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                              </w:rPr>
                            </w:pPr>
                            <w:r>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                                <w:rtl w:val="0"/>
                                <w:lang w:val="en-US"/>
                              </w:rPr>
                              <w:t xml:space="preserve">            throw new SaldoInsuficienteException();</w:t>
                            </w:r>
                          </w:p>
                          <w:p>
                            <w:pPr>
                              <w:pStyle w:val="Body"/>
                              <w:jc w:val="left"/>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                              </w:rPr>
                            </w:pPr>
                            <w:r>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                                <w:rtl w:val="0"/>
                              </w:rPr>
                              <w:t xml:space="preserve">        }</w:t>
                            </w:r>
                          </w:p>
                          <w:p>
                            <w:pPr>
                              <w:pStyle w:val="Body"/>
                              <w:jc w:val="left"/>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                              </w:rPr>
                            </w:pPr>
                            <w:r>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                                <w:rtl w:val="0"/>
                                <w:lang w:val="en-US"/>
                              </w:rPr>
                              <w:t xml:space="preserve">        String result = "*****************\n";</w:t>
                            </w:r>
                          </w:p>
                          <w:p>
                            <w:pPr>
                              <w:pStyle w:val="Body"/>
                              <w:jc w:val="left"/>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                              </w:rPr>
                            </w:pPr>
                            <w:r>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                                <w:rtl w:val="0"/>
                              </w:rPr>
                              <w:t xml:space="preserve">        result += "*** R$ " + saldo + ",00 ****\n";</w:t>
                            </w:r>
                          </w:p>
                          <w:p>
                            <w:pPr>
                              <w:pStyle w:val="Body"/>
                              <w:jc w:val="left"/>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                              </w:rPr>
                            </w:pPr>
                            <w:r>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                                <w:rtl w:val="0"/>
                                <w:lang w:val="en-US"/>
                              </w:rPr>
                              <w:t xml:space="preserve">        result += "*****************\n";</w:t>
                            </w:r>
                          </w:p>
                          <w:p>
                            <w:pPr>
                              <w:pStyle w:val="Body"/>
                              <w:jc w:val="left"/>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                              </w:rPr>
                            </w:pPr>
                            <w:r>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                                <w:rtl w:val="0"/>
                              </w:rPr>
                              <w:t xml:space="preserve">        return result;</w:t>
                            </w:r>
                          </w:p>
                          <w:p>
                            <w:pPr>
                              <w:pStyle w:val="Body"/>
                              <w:jc w:val="left"/>
                            </w:pPr>
                            <w:r>
                              <w:rPr>
                                <w:color w:val="fefefe"/>
                                <w:rtl w:val="0"/>
                              </w:rPr>
                              <w:t xml:space="preserve">    }</w:t>
                            </w:r>
                          </w:p>
                        </w:txbxContent>
                      </wps:txbx>
                      <wps:bodyPr wrap="square" lIns="50800" tIns="50800" rIns="50800" bIns="50800" numCol="1" anchor="t">
                        <a:noAutofit/>
                      </wps:bodyPr>
                    </wps:wsp>
                  </a:graphicData>
                </a:graphic>
              </wp:anchor>
            </w:drawing>
          </mc:Choice>
          <mc:Fallback>
            <w:pict>
              <v:rect id="_x0000_s1027" style="visibility:visible;position:absolute;margin-left:-0.5pt;margin-top:26.8pt;width:451.0pt;height:198.0pt;z-index:251660288;mso-position-horizontal:absolute;mso-position-horizontal-relative:margin;mso-position-vertical:absolute;mso-position-vertical-relative:line;mso-wrap-distance-left:12.0pt;mso-wrap-distance-top:12.0pt;mso-wrap-distance-right:12.0pt;mso-wrap-distance-bottom:12.0pt;">
                <v:fill color="#000000" opacity="100.0%" type="solid"/>
                <v:stroke on="f" weight="1.0pt" dashstyle="solid" endcap="flat" miterlimit="400.0%" joinstyle="miter" linestyle="single" startarrow="none" startarrowwidth="medium" startarrowlength="medium" endarrow="none" endarrowwidth="medium" endarrowlength="medium"/>
                <v:textbox>
                  <w:txbxContent>
                    <w:p>
                      <w:pPr>
                        <w:pStyle w:val="Body"/>
                        <w:jc w:val="left"/>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                        </w:rPr>
                      </w:pPr>
                      <w:r>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                          <w:rtl w:val="0"/>
                          <w:lang w:val="en-US"/>
                        </w:rPr>
                        <w:t xml:space="preserve"> public String imprimir() throws SaldoInsuficienteException {</w:t>
                      </w:r>
                    </w:p>
                    <w:p>
                      <w:pPr>
                        <w:pStyle w:val="Body"/>
                        <w:jc w:val="left"/>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                        </w:rPr>
                      </w:pPr>
                      <w:r>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                          <w:rtl w:val="0"/>
                        </w:rPr>
                        <w:t xml:space="preserve">        if (saldo &lt; valor) {</w:t>
                      </w:r>
                    </w:p>
                    <w:p>
                      <w:pPr>
                        <w:pStyle w:val="Body"/>
                        <w:jc w:val="left"/>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                        </w:rPr>
                      </w:pPr>
                      <w:r>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                          <w:rtl w:val="0"/>
                          <w:lang w:val="en-US"/>
                        </w:rPr>
                        <w:t xml:space="preserve">            throw new SaldoInsuficienteException();</w:t>
                      </w:r>
                    </w:p>
                    <w:p>
                      <w:pPr>
                        <w:pStyle w:val="Body"/>
                        <w:jc w:val="left"/>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                        </w:rPr>
                      </w:pPr>
                      <w:r>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                          <w:rtl w:val="0"/>
                        </w:rPr>
                        <w:t xml:space="preserve">        }</w:t>
                      </w:r>
                    </w:p>
                    <w:p>
                      <w:pPr>
                        <w:pStyle w:val="Body"/>
                        <w:jc w:val="left"/>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                        </w:rPr>
                      </w:pPr>
                      <w:r>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                          <w:rtl w:val="0"/>
                          <w:lang w:val="en-US"/>
                        </w:rPr>
                        <w:t xml:space="preserve">        String result = "*****************\n";</w:t>
                      </w:r>
                    </w:p>
                    <w:p>
                      <w:pPr>
                        <w:pStyle w:val="Body"/>
                        <w:jc w:val="left"/>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                        </w:rPr>
                      </w:pPr>
                      <w:r>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                          <w:rtl w:val="0"/>
                        </w:rPr>
                        <w:t xml:space="preserve">        result += "*** R$ " + saldo + ",00 ****\n";</w:t>
                      </w:r>
                    </w:p>
                    <w:p>
                      <w:pPr>
                        <w:pStyle w:val="Body"/>
                        <w:jc w:val="left"/>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                        </w:rPr>
                      </w:pPr>
                      <w:r>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                          <w:rtl w:val="0"/>
                          <w:lang w:val="en-US"/>
                        </w:rPr>
                        <w:t xml:space="preserve">        result += "*****************\n";</w:t>
                      </w:r>
                    </w:p>
                    <w:p>
                      <w:pPr>
                        <w:pStyle w:val="Body"/>
                        <w:jc w:val="left"/>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                        </w:rPr>
                      </w:pPr>
                      <w:r>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                          <w:rtl w:val="0"/>
                        </w:rPr>
                        <w:t xml:space="preserve">        return result;</w:t>
                      </w:r>
                    </w:p>
                    <w:p>
                      <w:pPr>
                        <w:pStyle w:val="Body"/>
                        <w:jc w:val="left"/>
                      </w:pPr>
                      <w:r>
                        <w:rPr>
                          <w:color w:val="fefefe"/>
                          <w:rtl w:val="0"/>
                        </w:rPr>
                        <w:t xml:space="preserve">    }</w:t>
                      </w:r>
                    </w:p>
                  </w:txbxContent>
                </v:textbox>
                <w10:wrap type="topAndBottom" side="bothSides" anchorx="margin"/>
              </v:rect>
            </w:pict>
          </mc:Fallback>
        </mc:AlternateContent>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Body"/>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>Para consertar esse erro basta acrescentar a linha de c</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>ó</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">digo na qual o valor do ticket </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t xml:space="preserve">é </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl w:val="0"/>
          <w:lang w:val="pt-PT"/>
        </w:rPr>
        <w:t>debitado do saldo total do cliente.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Body"/>
        <w:jc w:val="both"/>
      </w:pPr>
    </w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>'
$r.InsertXML($xml)

# --- Step 2: the document now spans a 2nd page, so the footer's cached PAGE field
#     result needs to go from "1" to "2" ---
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$ftrRange = $ftr.Range
$ftrRange.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "2", 2)
